$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cell A3 (previously "RUN") completely
$ws.Range("A3").ClearContents()

# Clear cell B3 value but keep formatting/style (previously "DGS-191")
$ws.Range("B3").ClearContents()

# Update the text in L2 from "...BNIMF 05" to "...BNIMF 08"
$ws.Range("L2").Value = "UFT Test Add Leads Prospek BNIMF 08"

# Change the top-left visible cell of the sheet view from J1 to F1
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Application.ActiveWindow.ScrollRow = 1
